$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# --- Insert new row 7: deadHeadTripAllowedModes ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "global"
$ws.Range("B7").Value = "deadHeadTripAllowedModes"
$ws.Range("C7").Value = "pt"
$ws.Range("D7").Clear()
$ws.Range("E7").Value = "Deadhead trips are routed using network links that match one of these allowed modes (list sperated with colon: pt,rail)"

# --- Insert new row 12: capacityFactor ---
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "global"
$ws.Range("B12").Value = "capacityFactor"
$ws.Range("C12").Value = 1.33
$ws.Range("D12").Clear()
$ws.Range("E12").Value = "Adjust the passenger capacity of units to reflect deviations in passenger demand."

# --- Refresh the autofilter range to cover the new rows ---
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:E26").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "scenario_info!_FilterDatabase") {
        $n.RefersTo = "=scenario_info!`$A`$1:`$E`$26"
    }
}

# --- Reflect final cursor position ---
$null = $ws.Range("C15").Select()
